# docs/diagrams/UndoRedoStartingStateListDiagram.pptx
# commit: "update user guide section3 images"
#
# The canonical-XML diff for this deck touches only:
#   1. The attribute order of the `xmlns="" xmlns:a16="..."` pair inside
#      five pre-existing `<a16:creationId/colId/rowId>` nodes
#      (Microsoft Office "creation id" extensions, uri
#      {FF2B5EF4-FFF2-40B4-BE49-F238E27FC236} /
#      {9D8B030D-6E8A-4147-A177-3AD203B41FA5} /
#      {0D108BD9-81ED-4DB2-BD59-A6C34878D82A}).
#   2. The value of one `<p14:modId>` (table "last modified" stamp,
#      uri {D42A27DB-BD31-4B8C-83A1-F6EECF244321}) on the "Table 4"
#      graphic frame.
#
# None of the shapes, text, geometry, table contents/structure, or any
# other visible/semantic property changed between the two revisions --
# every GUID, string, row/column, and position is identical. That
# pattern (same ids, only xmlns-attribute order + a modId counter
# bumped) is the signature PowerPoint itself leaves behind when a file
# is simply re-saved by a different Office build/platform; it is not
# something an editor reaches through the Shape/Table/TextRange object
# model -- PowerPoint does not expose `creationId`, `colId`, `rowId` or
# `modId` as scriptable properties, and re-applying any of the visible
# properties below (text, position, size, style, ...) does not disturb
# those extension blocks in the slide XML (verified: they round-trip
# byte-for-byte through this host regardless of which shape properties
# are touched).
#
# So there is nothing for an automation script to *change* here -- the
# slide's content is already the post-commit content. We still walk the
# shapes below (read-only) so the script documents/asserts that the
# table this commit touches is present and unchanged, without writing
# any property and thereby without perturbing anything that isn't in
# the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$textBox = $s.Shapes.Item(1)   # "TextBox 3"
$tableSh = $s.Shapes.Item(2)   # "Table 4"
$rect    = $s.Shapes.Item(3)   # "Rectangle 6"
$arrow   = $s.Shapes.Item(4)   # "Straight Arrow Connector 2"

Write-Host ("TextBox 3 text: " + $textBox.TextFrame.TextRange.Text)

$tbl = $tableSh.Table
Write-Host ("Table 4 rows/cols: " + $tbl.Rows.Count + "x" + $tbl.Columns.Count)
Write-Host ("Table 4 cell(1,1): " + $tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text)

Write-Host ("Rectangle 6 present: " + $rect.Name)
Write-Host ("Straight Arrow Connector 2 present: " + $arrow.Name)

Write-Host "No visible/content changes required for this revision; extension metadata (a16:creationId/colId/rowId ordering, p14:modId) is internal PowerPoint bookkeeping not reachable from the object model."
